$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1764
$ws1.Range("F5").Value = 443
$ws1.Range("F7").Value = 64
$ws1.Range("F9").Value = 316
$ws1.Range("F10").Value = 1703
$ws1.Range("F12").Value = 1407
$ws1.Range("F13").Value = 796
$ws1.Range("F15").Value = 668
$ws1.Range("F16").Value = 12690
$ws1.Range("F17").Value = 12710
$ws1.Range("F18").Value = 944
$ws1.Range("F21").Value = 502
$ws1.Range("F22").Value = 47
$ws1.Range("F23").Value = 522
$ws1.Range("F24").Value = 1990
$ws1.Range("F27").Value = 237

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 76
$ws2.Range("F10").Value = 69

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 159

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 159
$ws4.Range("F6").Value = 1764
$ws4.Range("F7").Value = 443
$ws4.Range("F10").Value = 64
$ws4.Range("F14").Value = 316
$ws4.Range("F15").Value = 1703
$ws4.Range("F17").Value = 1407
$ws4.Range("F18").Value = 796
$ws4.Range("F20").Value = 76
$ws4.Range("F21").Value = 668
$ws4.Range("F22").Value = 12690
$ws4.Range("F23").Value = 12710
$ws4.Range("F24").Value = 944
$ws4.Range("F27").Value = 502
$ws4.Range("F28").Value = 47
$ws4.Range("F29").Value = 522
$ws4.Range("F32").Value = 1990
$ws4.Range("F37").Value = 237
$ws4.Range("F39").Value = 69
